$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column D ("${regions}") — this shifts the old
#    D (Email), E (password), F (Description) columns one slot to the right
#    (-> E, F, G). Excel's column insert already carries forward the correct
#    per-cell styles (C's style into the new D, old D's style into E, etc.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# New column header
$ws.Range("D1").Value = "`${regions}"

# ---------------------------------------------------------------------------
# 2. Fix up the styles that an Excel column-insert does not carry over
#    correctly for the former *last* column (old F -> new G):
#      - G1 (header) must look like the other header cells (style of F1)
#      - G2:G9 (data rows) must look like column F (style of F2)
# ---------------------------------------------------------------------------
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("G2:G9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 2's old description comment is removed (G2 is now blank)
$ws.Range("G2").Value = ""

# ---------------------------------------------------------------------------
# 3. Data corrections
# ---------------------------------------------------------------------------
# Country code fix: UK -> GBR
$ws.Range("C5").Value = "GBR"

# Row 6: flag flipped to NO, and a new remark added in the description column
$ws.Range("B6").Value = "NO"
$ws.Range("G6").Value = "URL HS"

# Row 7 gets a region value
$ws.Range("D7").Value = "BEL"

# ---------------------------------------------------------------------------
# 4. Add two new test case rows (8 and 9), cloned from row 7's formatting
# ---------------------------------------------------------------------------
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)
$ws.Range("A9:G9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A8").Value = "UC7_SFCC_Connect_UAT"
$ws.Range("B8").Value = "YES"
$ws.Range("C8").Value = "BNL"
$ws.Range("D8").Value = "LUX"
$ws.Range("E8").Value = "UC6_UAT@gg.com"
$ws.Range("F8").Value = "azertyui"
$ws.Range("G8").Value = ""

$ws.Range("A9").Value = "UC8_SFCC_Connect_UAT"
$ws.Range("B9").Value = "YES"
$ws.Range("C9").Value = "BNL"
$ws.Range("D9").Value = "NLD"
$ws.Range("E9").Value = "UC6_UAT@gg.com"
$ws.Range("F9").Value = "azertyui"
$ws.Range("G9").Value = ""

# ---------------------------------------------------------------------------
# 5. Update conditional formatting ranges to cover the new layout
# ---------------------------------------------------------------------------
$fcsB = $ws.Range("B2").FormatConditions
$fcsB.Item(1).ModifyAppliesToRange($ws.Range("B2:B9"))
$fcsB.Item(2).ModifyAppliesToRange($ws.Range("B2:B9"))

$fcsC = $ws.Range("C2").FormatConditions
$fcsC.Item(1).ModifyAppliesToRange($ws.Range("C2:D2"))
$fcsC.Item(2).ModifyAppliesToRange($ws.Range("C2:D2"))

$fcsE = $ws.Range("D1").FormatConditions
$fcsE.Item(1).ModifyAppliesToRange($ws.Range("E1"))
$fcsE.Item(2).ModifyAppliesToRange($ws.Range("E1"))

# ---------------------------------------------------------------------------
# 6. Restore selection as left by the author
# ---------------------------------------------------------------------------
$ws.Range("B7").Select()
